$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 2842.8647
$ws.Range("I15").Value = 2842.8647
$ws.Range("K15").Value = 8528.5941
$ws.Range("M15").Value = -8359.5941
# row 92
$ws.Range("H92").Value = 1913.3334
$ws.Range("I92").Value = 2012.8572
$ws.Range("K92").Value = 2012.8572
$ws.Range("M92").Value = -764.8571999999999
# row 98
$ws.Range("H98").Value = 2761.6843
$ws.Range("J98").Value = 1899
$ws.Range("L98").Value = 1899
$ws.Range("N98").Value = -4895
# row 122
$ws.Range("H122").Value = 2761.6843
$ws.Range("J122").Value = 1899
$ws.Range("L122").Value = 5697
$ws.Range("N122").Value = -10597
# row 132
$ws.Range("H132").Value = 6672370
$ws.Range("I132").Value = 10106597
$ws.Range("J132").Value = 5929.5293
$ws.Range("K132").Value = 30319791
$ws.Range("L132").Value = 17788.5879
$ws.Range("M132").Value = -30317261
$ws.Range("N132").Value = -22848.5879
# row 137
$ws.Range("H137").Value = 1121.8939
$ws.Range("I137").Value = 805.3415
$ws.Range("J137").Value = 1641.04
$ws.Range("K137").Value = 2416.0245
$ws.Range("L137").Value = 4923.12
$ws.Range("M137").Value = 133.9755
$ws.Range("N137").Value = -10023.12
# row 138
$ws.Range("H138").Value = 1330.73
$ws.Range("I138").Value = 745.25
$ws.Range("J138").Value = 1606.25
$ws.Range("K138").Value = 2235.75
$ws.Range("L138").Value = 4818.75
$ws.Range("M138").Value = 2904.25
$ws.Range("N138").Value = -15098.75
# row 141
$ws.Range("H141").Value = 583.0476
$ws.Range("I141").Value = 512.2
$ws.Range("K141").Value = 1536.6
$ws.Range("M141").Value = 3643.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4259.655
$ws.Range("I32").Value = 3709.7778
$ws.Range("K32").Value = 3709.7778
$ws.Range("M32").Value = -3422.7778
# row 62
$ws.Range("H62").Value = 60000
$ws.Range("J62").Value = 60000
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61248
# row 63
$ws.Range("H63").Value = 25001980
$ws.Range("I63").Value = 1853.1111
$ws.Range("K63").Value = 1853.1111
$ws.Range("M63").Value = -1167.1111
# row 65
$ws.Range("H65").Value = 60000
$ws.Range("J65").Value = 60000
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186240
# row 66
$ws.Range("H66").Value = 25001980
$ws.Range("I66").Value = 1853.1111
$ws.Range("K66").Value = 9265.5555
$ws.Range("M66").Value = -5833.5555
# row 110
$ws.Range("H110").Value = 1295.3214
$ws.Range("I110").Value = 876
$ws.Range("J110").Value = 1714.6428
$ws.Range("K110").Value = 876
$ws.Range("L110").Value = 1714.6428
$ws.Range("M110").Value = 1169
$ws.Range("N110").Value = -5804.6428
# row 132
$ws.Range("H132").Value = 2424
$ws.Range("I132").Value = 2277.5833
$ws.Range("J132").Value = 2599.7
$ws.Range("K132").Value = 6832.749899999999
$ws.Range("L132").Value = 7799.099999999999
$ws.Range("M132").Value = -4302.749899999999
$ws.Range("N132").Value = -12859.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 3833.275
$ws.Range("I134").Value = 1126.5161
$ws.Range("J134").Value = 13156.556
$ws.Range("K134").Value = 3379.5483
$ws.Range("L134").Value = 39469.66800000001
$ws.Range("M134").Value = -844.5483000000004
$ws.Range("N134").Value = -44539.66800000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2351.5789
$ws.Range("I31").Value = 2510.5881
$ws.Range("K31").Value = 2510.5881
$ws.Range("M31").Value = -2215.5881
# row 34
$ws.Range("H34").Value = 2351.5789
$ws.Range("I34").Value = 2510.5881
$ws.Range("K34").Value = 2510.5881
$ws.Range("M34").Value = -2308.5881
# row 58
$ws.Range("H58").Value = 1133.1389
$ws.Range("I58").Value = 1082.1666
$ws.Range("J58").Value = 1388
$ws.Range("K58").Value = 1082.1666
$ws.Range("L58").Value = 1388
$ws.Range("M58").Value = -879.1666
$ws.Range("N58").Value = -1794
# row 132
$ws.Range("H132").Value = 1987.4546
$ws.Range("I132").Value = 1645
$ws.Range("J132").Value = 2672.3635
$ws.Range("K132").Value = 4935
$ws.Range("L132").Value = 8017.0905
$ws.Range("M132").Value = -2405
$ws.Range("N132").Value = -13077.0905
# row 134
$ws.Range("H134").Value = 15626148
$ws.Range("I134").Value = 1086.5714
$ws.Range("J134").Value = 45455812
$ws.Range("K134").Value = 3259.7142
$ws.Range("L134").Value = 136367436
$ws.Range("M134").Value = -724.7142000000003
$ws.Range("N134").Value = -136372506
# row 136
$ws.Range("H136").Value = 1133.1389
$ws.Range("I136").Value = 1082.1666
$ws.Range("J136").Value = 1388
$ws.Range("K136").Value = 3246.4998
$ws.Range("L136").Value = 4164
$ws.Range("M136").Value = -696.4998000000001
$ws.Range("N136").Value = -9264

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 112.375
$ws.Range("I2").Value = 66.333336
$ws.Range("J2").Value = 140
$ws.Range("K2").Value = 398.000016
$ws.Range("L2").Value = 840
$ws.Range("M2").Value = -285.000016
$ws.Range("N2").Value = -1066
# row 22
$ws.Range("H22").Value = 2000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 6000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -6338
# row 27
$ws.Range("H27").Value = 2000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 6000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -6204
# row 56
$ws.Range("H56").Value = 5728.9414
$ws.Range("I56").Value = 5728.9414
$ws.Range("K56").Value = 5728.9414
$ws.Range("M56").Value = -5198.9414
# row 131
$ws.Range("H131").Value = 25644598
$ws.Range("J131").Value = 4364.355
$ws.Range("L131").Value = 13093.065
$ws.Range("N131").Value = -23173.065
# row 132
$ws.Range("H132").Value = 1484.2106
$ws.Range("I132").Value = 1053.8462
$ws.Range("J132").Value = 2416.6667
$ws.Range("K132").Value = 9484.6158
$ws.Range("L132").Value = 21750.0003
$ws.Range("M132").Value = -6954.6158
$ws.Range("N132").Value = -26810.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 40911790
$ws.Range("I70").Value = 41669316
$ws.Range("J70").Value = 40002760
$ws.Range("K70").Value = 41669316
$ws.Range("L70").Value = 40002760
$ws.Range("M70").Value = -41669046
$ws.Range("N70").Value = -40003300
# row 73
$ws.Range("H73").Value = 40911790
$ws.Range("I73").Value = 41669316
$ws.Range("J73").Value = 40002760
$ws.Range("K73").Value = 41669316
$ws.Range("L73").Value = 40002760
$ws.Range("M73").Value = -41668380
$ws.Range("N73").Value = -40004632
# row 107
$ws.Range("H107").Value = 712838.9399999999
$ws.Range("I107").Value = 1603439.1
$ws.Range("J107").Value = 358.86667
$ws.Range("K107").Value = 1603439.1
$ws.Range("L107").Value = 358.86667
$ws.Range("M107").Value = -1601519.1
$ws.Range("N107").Value = -4198.86667
# row 126
$ws.Range("H126").Value = 1669.4
$ws.Range("I126").Value = 1349.25
$ws.Range("K126").Value = 4047.75
$ws.Range("M126").Value = -1577.75
# row 132
$ws.Range("H132").Value = 2945.6667
$ws.Range("I132").Value = 2625.5
$ws.Range("J132").Value = 3411.3635
$ws.Range("K132").Value = 7876.5
$ws.Range("L132").Value = 10234.0905
$ws.Range("M132").Value = -5346.5
$ws.Range("N132").Value = -15294.0905
# row 136
$ws.Range("H136").Value = 25990.727
$ws.Range("J136").Value = 25990.727
$ws.Range("L136").Value = 77972.181
$ws.Range("N136").Value = -83072.181

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 93
$ws.Range("H93").Value = 582.1
$ws.Range("I93").Value = 536.125
$ws.Range("J93").Value = 766
$ws.Range("K93").Value = 536.125
$ws.Range("L93").Value = 766
$ws.Range("M93").Value = 711.875
$ws.Range("N93").Value = -3262
# row 136
$ws.Range("H136").Value = 2501.25
$ws.Range("I136").Value = 2500
$ws.Range("J136").Value = 2501.6667
$ws.Range("K136").Value = 7500
$ws.Range("L136").Value = 7505.000100000001
$ws.Range("M136").Value = -4950
$ws.Range("N136").Value = -12605.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 82
$ws.Range("H82").Value = 13500
$ws.Range("J82").Value = 13500
$ws.Range("L82").Value = 13500
$ws.Range("N82").Value = -14266
# row 85
$ws.Range("H85").Value = 13500
$ws.Range("J85").Value = 13500
$ws.Range("L85").Value = 13500
$ws.Range("N85").Value = -16152
# row 113
$ws.Range("H113").Value = 399.9524
$ws.Range("I113").Value = 351.375
$ws.Range("J113").Value = 429.84616
$ws.Range("K113").Value = 1054.125
$ws.Range("L113").Value = 1289.53848
$ws.Range("M113").Value = 1115.875
$ws.Range("N113").Value = -5629.53848
# row 132
$ws.Range("H132").Value = 1599.6327
$ws.Range("I132").Value = 1089.4839
$ws.Range("K132").Value = 3268.4517
$ws.Range("M132").Value = -738.4516999999996
# row 136
$ws.Range("H136").Value = 899.55
$ws.Range("I136").Value = 867.9286
$ws.Range("J136").Value = 973.3333
$ws.Range("K136").Value = 2603.7858
$ws.Range("L136").Value = 2919.9999
$ws.Range("M136").Value = -53.78579999999965
$ws.Range("N136").Value = -8019.9999
# row 137
$ws.Range("H137").Value = 30649.9
$ws.Range("J137").Value = 30649.9
$ws.Range("L137").Value = 30649.9
$ws.Range("N137").Value = -40849.9
# row 139
$ws.Range("H139").Value = 34647.5
$ws.Range("J139").Value = 34647.5
$ws.Range("L139").Value = 34647.5
$ws.Range("N139").Value = -44927.5
